$d = $word.ActiveDocument

function Set-ParagraphText($paraIndex, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $newText
}

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Output ("WARNING: text not found in paragraph " + $paraIndex)
    }
}

Set-ParagraphText 6 "Diagramas para estudos de processos químicos. 2  Estrutura e síntese de processos químicos industriais. 3  Análise de desempenho de processos químicos. 4  Estudo de planta química industrial"
Set-ParagraphText 7 "1 - Diagrams for Understanding Chemical Processes. 2 - The Structure and Synthesis of Process Flow Diagrams. 3 - Analysis of process performance. 4 - Industrial chemical plant study"
Set-ParagraphText 9 "1 - Consolidação e aplicação dos conhecimentos adquiridos em cada uma das áreas específicas do curso de Engenharia Química.  2  Integração dos conhecimentos da Engenharia Química."
Set-ParagraphText 11 "1 - Diagramas para estudos de processos químicos: diagramas de bloco; Fluxogramas de processo (PFD); Fluxogramas de instrumentação e tubulação (P&ID).2  Estrutura e síntese de processos químicos industriais: Hierarquia no planejamento de processos; Etapa 1- Descontínuo ou contínuo; Etapa 2 - Estrutura de entrada/saída de processo; Etapa 3- Estrutura de reciclo; 3  Análise de desempenho de processos químicos: Modelo de entrada e saída; Ferramentas para a avaliação de processos.4  Estudo de planta química industrial."
Set-ParagraphText 12 "1 - Consolidation and application of knowledge acquired in each of the specific areas of the Chemical Engineering degree. 2 - Integration of knowledge of Chemical Engineering"
Set-ParagraphText 14 "Provas escritas e Apresentação de Trabalhos"
Set-ParagraphText 19 "5816812 - João Paulo Alves Silva"

# Paragraph 17 has 3 bold labels interleaved with content runs; only content runs change.
# Process in reverse (last content run first) to avoid newly-written text colliding with next search.
Replace-InParagraph 17 "Média Final = (N + Prova Recuperação)/2" "PERLINGEIRO, Carlos A. G. Engenharia de processos: análise, simulação, otimização e síntese de processos químicos.  Editora Blucher, 2005.TURTON, BAILIE; WHITING; SHAEIWITZ  Analysis, Synthesis, and Design of Chemical Processes. 3. Ed. LTC Editora, 2005.COULSON, J. M.; RICHARDSON, J.F. Chemical Engineering Design: Chemical Engineering Volume 6. Editora Fourth, 2005.HIMMELBLAU, David M. Engenharia química princípios e cálculos. LTC Editora, 2006.FELDER, R.M; Rousseau, R.W. Princípios elementares dos processos químicos. LTC Editora, 2005."
Replace-InParagraph 17 "A nota será composta por ao menos uma prova escrita e trabalhos realizados e apresentados durante o semestre. O peso de cada atividade será definido segundo critérios do professor." "Média Final = (N + Prova Recuperação)/2"
Replace-InParagraph 17 "Provas escritas e Apresentação de Trabalhos" "A nota será composta por ao menos uma prova escrita e trabalhos realizados e apresentados durante o semestre. O peso de cada atividade será definido segundo critérios do professor."
